$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INPUT_SHEET")

# Update existing rows 2-6 (BAN values change, tickets keep same column position
# but their text content changes)
$ws.Cells.Item(2,1).Value = 155283019
$ws.Cells.Item(2,2).Value = "ONREG-19630"

$ws.Cells.Item(3,1).Value = 162691762
$ws.Cells.Item(3,2).Value = "ONREG-19634"

$ws.Cells.Item(4,1).Value = 202592941
$ws.Cells.Item(4,2).Value = "ONREG-19705"

$ws.Cells.Item(5,1).Value = 288561502
$ws.Cells.Item(5,2).Value = "ONREG-25688"

$ws.Cells.Item(6,1).Value = 266720569
$ws.Cells.Item(6,2).Value = "ONREG-20095"

# Rows 7-10 get new values (row 7 & 8 reuse the same BAN/ticket combos as
# the new rows 2 and 4 respectively)
$ws.Cells.Item(7,1).Value = 155283019
$ws.Cells.Item(7,2).Value = "ONREG-19630"

$ws.Cells.Item(8,1).Value = 202592941
$ws.Cells.Item(8,2).Value = "ONREG-19705"

$ws.Cells.Item(9,1).Value = 103127851
$ws.Cells.Item(9,2).Value = "ONREG-19939"

$ws.Cells.Item(10,1).Value = 203320557
$ws.Cells.Item(10,2).Value = "ONREG-20001"

# New rows 11-13: copy formatting from row 10 first (column style defaults to
# Text, so formats must be pasted before the value is written or the value
# gets stored as a string instead of a number), then write the values.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11,1).Value = 135391535
$ws.Cells.Item(11,2).Value = "ONREG-18620"

$ws.Cells.Item(12,1).Value = 121245443
$ws.Cells.Item(12,2).Value = "ONREG-26645"

$ws.Cells.Item(13,1).Value = 131507015
$ws.Cells.Item(13,2).Value = "ST1-27775"

# Selection / active-view bookkeeping: INPUT_SHEET becomes the active tab and
# its selection moves to C3; API_CONSOLE_INFO loses tabSelected as a result.
$ws.Activate()
$ws.Range("C3").Select()
